$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 787.0769
$ws.Range("J17").Value = 751.087
$ws.Range("L17").Value = 2253.261
$ws.Range("N17").Value = -2589.261
$ws.Range("H40").Value = 1699.1111
$ws.Range("I40").Value = 1144.4
$ws.Range("J40").Value = 2392.5
$ws.Range("K40").Value = 1144.4
$ws.Range("L40").Value = 2392.5
$ws.Range("M40").Value = -969.4000000000001
$ws.Range("N40").Value = -2742.5
$ws.Range("H116").Value = 33083
$ws.Range("I116").Value = 25724.5
$ws.Range("J116").Value = 35185.43
$ws.Range("K116").Value = 25724.5
$ws.Range("L116").Value = 35185.43
$ws.Range("M116").Value = -22282.5
$ws.Range("N116").Value = -42069.43
$ws.Range("H132").Value = 20843
$ws.Range("I132").Value = 29263.428
$ws.Range("K132").Value = 87790.284
$ws.Range("M132").Value = -85260.284
$ws.Range("H137").Value = 3818.8276
$ws.Range("I137").Value = 2028.5834
$ws.Range("J137").Value = 12412
$ws.Range("K137").Value = 6085.7502
$ws.Range("L137").Value = 37236
$ws.Range("M137").Value = -3535.7502
$ws.Range("N137").Value = -42336
$ws.Range("H138").Value = 3975.6943
$ws.Range("I138").Value = 3513.6206
$ws.Range("K138").Value = 10540.8618
$ws.Range("M138").Value = -5400.861800000001
$ws.Range("H141").Value = 5019.5
$ws.Range("I141").Value = 2595.6956
$ws.Range("J141").Value = 12983.429
$ws.Range("K141").Value = 7787.0868
$ws.Range("L141").Value = 38950.287
$ws.Range("M141").Value = -2607.0868
$ws.Range("N141").Value = -49310.287

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1077.4
$ws.Range("I2").Value = 1048
$ws.Range("K2").Value = 1048
$ws.Range("M2").Value = -935
$ws.Range("H74").Value = 6618.64
$ws.Range("I74").Value = 3653.634
$ws.Range("K74").Value = 3653.634
$ws.Range("M74").Value = -2779.634
$ws.Range("H77").Value = 6618.64
$ws.Range("I77").Value = 3653.634
$ws.Range("K77").Value = 18268.17
$ws.Range("M77").Value = -13900.17
$ws.Range("H110").Value = 1945
$ws.Range("I110").Value = 1899.5
$ws.Range("K110").Value = 1899.5
$ws.Range("M110").Value = 145.5
$ws.Range("H116").Value = 1077.4
$ws.Range("I116").Value = 1048
$ws.Range("K116").Value = 1048
$ws.Range("M116").Value = 1246
$ws.Range("H132").Value = 8328.25
$ws.Range("I132").Value = 5964.5454
$ws.Range("K132").Value = 17893.6362
$ws.Range("M132").Value = -15363.6362

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1077.4
$ws.Range("I3").Value = 1048
$ws.Range("K3").Value = 1048
$ws.Range("M3").Value = -934
$ws.Range("H75").Value = 21038
$ws.Range("J75").Value = 55000
$ws.Range("L75").Value = 55000
$ws.Range("N75").Value = -56872
$ws.Range("H78").Value = 21038
$ws.Range("J78").Value = 55000
$ws.Range("L78").Value = 165000
$ws.Range("N78").Value = -174360
$ws.Range("H80").Value = 556.7083
$ws.Range("J80").Value = 379.73334
$ws.Range("L80").Value = 379.73334
$ws.Range("N80").Value = -2375.73334
$ws.Range("H83").Value = 556.7083
$ws.Range("J83").Value = 379.73334
$ws.Range("L83").Value = 1898.6667
$ws.Range("N83").Value = -11882.6667
$ws.Range("H94").Value = 3325
$ws.Range("I94").Value = 3128.9688
$ws.Range("K94").Value = 3128.9688
$ws.Range("M94").Value = -2677.9688

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3340.6155
$ws.Range("I31").Value = 3143.2
$ws.Range("J31").Value = 3998.6667
$ws.Range("K31").Value = 3143.2
$ws.Range("L31").Value = 3998.6667
$ws.Range("M31").Value = -2848.2
$ws.Range("N31").Value = -4588.6667
$ws.Range("H34").Value = 3340.6155
$ws.Range("I34").Value = 3143.2
$ws.Range("J34").Value = 3998.6667
$ws.Range("K34").Value = 3143.2
$ws.Range("L34").Value = 3998.6667
$ws.Range("M34").Value = -2941.2
$ws.Range("N34").Value = -4402.6667
$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("N50").Value = -51250
$ws.Range("H105").Value = 7269.353
$ws.Range("I105").Value = 7567.4375
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 7567.4375
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -5820.4375
$ws.Range("N105").Value = -5994

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 201.7
$ws.Range("J23").Value = 225.14285
$ws.Range("L23").Value = 675.4285500000001
$ws.Range("N23").Value = -1145.42855
$ws.Range("H87").Value = 26166.834
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 26166.834
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 78500.50199999999
$ws.Range("N87").Value = -80996.50199999999
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 26166.834
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 26166.834
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 235501.506
$ws.Range("N90").Value = -247981.506
$ws.Range("M90").ClearContents()
$ws.Range("H101").Value = 14298286
$ws.Range("J101").Value = 14298286
$ws.Range("L101").Value = 42894858
$ws.Range("N101").Value = -42899726

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 16000
$ws.Range("I18").Value = 16000
$ws.Range("K18").Value = 16000
$ws.Range("M18").Value = -15707
$ws.Range("H21").Value = 6098.125
$ws.Range("I21").Value = 4827.143
$ws.Range("K21").Value = 4827.143
$ws.Range("M21").Value = -4654.143
$ws.Range("H30").Value = 6098.125
$ws.Range("I30").Value = 4827.143
$ws.Range("K30").Value = 4827.143
$ws.Range("M30").Value = -4722.143
$ws.Range("H80").Value = 1688.75
$ws.Range("J80").Value = 1378
$ws.Range("L80").Value = 1378
$ws.Range("N80").Value = -3374
$ws.Range("H83").Value = 1688.75
$ws.Range("J83").Value = 1378
$ws.Range("L83").Value = 6890
$ws.Range("N83").Value = -16874
$ws.Range("H122").Value = 3687.0476
$ws.Range("I122").Value = 3702.3572
$ws.Range("J122").Value = 3656.4285
$ws.Range("K122").Value = 11107.0716
$ws.Range("L122").Value = 10969.2855
$ws.Range("M122").Value = -8657.071599999999
$ws.Range("N122").Value = -15869.2855

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6000
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888
$ws.Range("H13").Value = 8200
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 16000
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = -260
$ws.Range("N13").Value = -16280
$ws.Range("H46").Value = 3991.6943
$ws.Range("J46").Value = 4719.28
$ws.Range("L46").Value = 4719.28
$ws.Range("N46").Value = -5095.28
$ws.Range("H122").Value = 4052.5
$ws.Range("I122").Value = 2841.3333
$ws.Range("K122").Value = 8523.999899999999
$ws.Range("M122").Value = -6073.999899999999
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("H132").Value = 2598.1304
$ws.Range("I132").Value = 1815.1177
$ws.Range("J132").Value = 4816.6665
$ws.Range("K132").Value = 5445.3531
$ws.Range("L132").Value = 14449.9995
$ws.Range("M132").Value = -2915.3531
$ws.Range("N132").Value = -19509.9995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2750.8333
$ws.Range("I126").Value = 2300
$ws.Range("K126").Value = 6900
$ws.Range("M126").Value = -4430
